$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column W (after the existing last date column V) and carry
# over that column's formatting (bold header, thin border, centered) so the
# new header cell is visually consistent with the other date headers.
$ws.Columns.Item(23).Insert()

# Header: new report date. Force the cell to remain literal text (like the
# other date headers, e.g. "07-10-2020" in V1) instead of Excel silently
# converting the recognizable date string into a date serial number.
$ws.Range("W1").NumberFormat = "@"
$ws.Range("W1").Value = "08-10-2020"

# Body: updated "active cases" counts for 08-10-2020, one per state/UT row
# (rows 2-36), mirroring the order already used in column V.
$values = @(
  185, 49513, 2850, 31786, 11326, 1448, 26777, 108, 22186, 4749,
  16485, 11029, 2996, 12131, 9759, 116172, 92246, 1228, 17522, 244976,
  2805, 2411, 231, 1200, 26368, 4680, 11563, 21351, 570, 45135,
  26368, 4389, 8367, 43154, 28361
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 23).Value = $values[$i]
}
